# "results for combined filters experiment"
#
# The CNN results log had an extra "Test Accuracy (shhhh)" column (I) that
# was a pending/placeholder metric. Remove it so the log only reports the
# combined-filters-experiment columns, shifting the trailing note column
# (old J, "<-low train accuracy = train longer") left into I.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CNN")

# Make sure we're looking at the right sheet (it was already the active one).
$ws.Activate()

# Select column I (the "Test Accuracy (shhhh)" column) the way a user would
# before deleting it, then delete the whole column - this shifts everything
# to its right (the stray note in column J) one place to the left.
$ws.Columns.Item(9).Select()
$ws.Columns.Item(9).EntireColumn.Delete()
